# Apply cryptos list price/volume updates for Sun Dec 3 02:43:47 UTC 2023 run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.450.58'
$ws.Range("E2").Value = '  +1.94%  '
$ws.Range("D3").Value = '2.170.07'
$ws.Range("E3").Value = '  +3.75%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.23'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +0.48%  '
$ws.Range("E6").Value = '  +1.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '65.29'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +6.92%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.401'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +4.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0866'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +2.80%  '
$ws.Range("E11").Value = '  +0.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.02'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +5.11%  '
$ws.Range("D13").Value = '2.492.01'
$ws.Range("E13").Value = '  +3.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.53'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +2.28%  '
$ws.Range("E15").Value = '  +0.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.58'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +1.89%  '
$ws.Range("D17").Value = '2.177.87'
$ws.Range("E17").Value = '  +4.27%  '
$ws.Range("D18").Value = '39.444.75'
$ws.Range("E18").Value = '  +2.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.20'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +1.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.38'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("E21").Value = '  +2.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '232.41'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +2.20%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.36'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -1.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +1.80%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.73'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +1.65%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '172.69'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +0.97%  '
$ws.Range("E28").Value = '  -0.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.15'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +4.71%  '
$ws.Range("E30").Value = '  -1.35%  '
$ws.Range("E31").Value = '  +12.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.123'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +2.03%  '
$ws.Range("E33").Value = '  +3.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.82'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +2.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.16'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +8.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0620'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +1.69%  '
$ws.Range("E37").Value = '  +1.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.60'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '104.89'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +3.96%  '
$ws.Range("E41").Value = '  +0.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.92'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -0.33%  '
$ws.Range("D43").Value = '1.541.32'
$ws.Range("E43").Value = '  +0.51%  '
$ws.Range("E44").Value = '  +5.47%  '
$ws.Range("E45").Value = '  +6.29%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.96'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +2.32%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.11'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +7.87%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0928'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  +1.23%  '
$ws.Range("B49").Value = 'HuobiToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.82'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +0.49%  '
$ws.Range("D50").Value = '2.375.00'
$ws.Range("E50").Value = '  +3.75%  '
$ws.Range("E51").Value = '  +0.30%  '
